# Auto update: 2025-12-03 03:05:24
# Updates the DECISION sheet with fresh data for 2025-12-03.
# Row 2 becomes UnitedHealth (UNH), Row 3 becomes Prudential (PRU) -- the two
# swap places -- and all metric columns are refreshed for every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The date column (A) holds text like "2025-12-03" rather than a real date
# serial. Excel auto-converts a bare date-shaped string typed into a cell,
# so force Text formatting first, write the value, then drop back to the
# sheet's default (General) formatting so no stray number format lingers
# on the cell.
$ws.Range("A2:A5").NumberFormat = "@"

# Row 2: UnitedHealth Group Incorporated (UNH)
$ws.Range("A2").Value = "2025-12-03"
$ws.Range("B2").Value = "UnitedHealth Group Incorporated"
$ws.Range("C2").Value = "UNH"
$ws.Range("D2").Value = 323.88
$ws.Range("E2").Value = 47.5
$ws.Range("F2").Value = 1.51
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 50
$ws.Range("I2").Value = 46
$ws.Range("J2").Value = 46
$ws.Range("K2").Value = 56
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 65.32892478746797
$ws.Range("O2").Value = "🟢 상승 우위 (다소 완화)"

# Row 3: Prudential Financial, Inc. (PRU)
$ws.Range("A3").Value = "2025-12-03"
$ws.Range("B3").Value = "Prudential Financial, Inc."
$ws.Range("C3").Value = "PRU"
$ws.Range("D3").Value = 107.68
$ws.Range("E3").Value = 59.1
$ws.Range("F3").Value = 0.46
$ws.Range("G3").Value = 60
$ws.Range("H3").Value = 46
$ws.Range("I3").Value = 46
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 56
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 65.32892478746797
$ws.Range("O3").Value = "🟢 상승 우위 (다소 완화)"

# Row 4: MetLife, Inc. (MET)
$ws.Range("A4").Value = "2025-12-03"
$ws.Range("B4").Value = "MetLife, Inc."
$ws.Range("C4").Value = "MET"
$ws.Range("D4").Value = 76.05
$ws.Range("E4").Value = 42.7
$ws.Range("F4").Value = 1.14
$ws.Range("G4").Value = 40
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 33
$ws.Range("K4").Value = 51.6
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 65.32892478746797
$ws.Range("O4").Value = "🟢 상승 우위 (다소 완화)"

# Row 5: American International Group, I (AIG)
$ws.Range("A5").Value = "2025-12-03"
$ws.Range("B5").Value = "American International Group, I"
$ws.Range("C5").Value = "AIG"
$ws.Range("D5").Value = 76.61
$ws.Range("E5").Value = 52.6
$ws.Range("F5").Value = 1.99
$ws.Range("G5").Value = 40
$ws.Range("H5").Value = 36
$ws.Range("I5").Value = 46
$ws.Range("J5").Value = 46
$ws.Range("K5").Value = 50
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 65.32892478746797
$ws.Range("O5").Value = "🟢 상승 우위 (다소 완화)"

# Drop the temporary Text format so the date cells end up back on General,
# matching the rest of the sheet's (unstyled) data rows.
$ws.Range("A2:A5").ClearFormats()
